# Commit: "Added new screens to split content"
#
# On the "Translation" sheet, rows 9-12 hold a block of "SingleUseId#" rows.
# A new screen/text ("SingleUseId13") was introduced, shifting the Text ID
# (column B) of each existing row down one slot in the series, while the
# English ("GB", column E) text content rotates up by one row, with the
# value that was on row 9 wrapping around onto the newly appended row 12.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

# Column B (Text ID) - renumber the SingleUseId series
$ws.Range("B9").Value  = "SingleUseId10"
$ws.Range("B10").Value = "SingleUseId11"
$ws.Range("B11").Value = "SingleUseId12"
$ws.Range("B12").Value = "SingleUseId13"

# Column E (GB / English text) - rotate values up by one row, wrapping around
$ws.Range("E9").Value  = "Send TCP"
$ws.Range("E10").Value = "Log Data"
$ws.Range("E11").Value = "Dump Log"
$ws.Range("E12").Value = "Toggle LED"
